$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Experiment Results")

$ws.Range("A10").Value = "AGGREGATION"
$ws.Range("A22").Value = "LARGE"
$ws.Range("A26").Value = "COUNT(*)"
$ws.Range("B26").Value = "manual_test_agg_all_3"
$ws.Range("B27").Value = "auto_test_agg_all_3"
$ws.Range("B28").Value = "manual_test_agg_all_5"
$ws.Range("B29").Value = "auto_test_agg_all_5"
